# Snapshot/release split: the demo-open step now reports its own
# pass/fail outcome in the "Resultat" column (instead of leaving it
# blank). "Succes" / "Echec" are the two possible outcomes; this run
# failed to open the demo, so every row is flagged with the failure
# message, shown in the same "ok" green used elsewhere in the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$failureText = "Échec : Ouverture DEMO."

$resultRange = $ws.Range("F2:F6")
$resultRange.Value = $failureText
$resultRange.Font.ColorIndex = 10
